$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# Row 8 data: move from 4th quarter 2022 reporting period to 2nd quarter 2023
$ws.Cells.Item(8, 1).Value = 2023
$ws.Cells.Item(8, 2).Value = 44927
$ws.Cells.Item(8, 3).Value = 45107
$ws.Cells.Item(8, 33).Value = 45117
$ws.Cells.Item(8, 34).Value = 45117

# AF8/AI8 text updates (shared strings)
$ws.Cells.Item(8, 32).Value = "Departamento de Infraestructura (UPP) "
$ws.Cells.Item(8, 35).Value = "La Universidad Politécnica de Pachuca, no cuenta con inventario de bienes inmuebles para el periodo reportado. "

$ws.Rows.Item(8).RowHeight = 60
